# Update "想去人数" (number of people interested) counts on the
# "展览" (Exhibitions) and "全部类型" (All Types) worksheets.
# These two sheets mirror the same underlying event data, so the same
# set of updates (by name match would also work) is applied to both,
# using the row offsets appropriate to each sheet.

$wb = $excel.ActiveWorkbook

$sheetExhibition = $wb.Worksheets.Item("展览")
$sheetAllTypes   = $wb.Worksheets.Item("全部类型")

# Row -> new value updates for the "展览" sheet (column F)
$updatesExhibition = @{
    6  = 194
    7  = 4542
    14 = 178
    15 = 963
    16 = 74
    20 = 108
    22 = 3455
    23 = 5794
    29 = 3345
    34 = 518
    36 = 204
    37 = 259
    38 = 347
    39 = 119
    41 = 898
    45 = 42
    47 = 61
}

foreach ($row in $updatesExhibition.Keys) {
    $sheetExhibition.Range("F$row").Value = $updatesExhibition[$row]
}

# Row -> new value updates for the "全部类型" sheet (column F)
$updatesAllTypes = @{
    6  = 194
    7  = 4542
    15 = 178
    16 = 963
    17 = 74
    21 = 108
    23 = 3455
    24 = 5794
    30 = 3345
    35 = 518
    37 = 204
    38 = 259
    39 = 347
    40 = 119
    42 = 898
    46 = 42
    48 = 61
}

foreach ($row in $updatesAllTypes.Keys) {
    $sheetAllTypes.Range("F$row").Value = $updatesAllTypes[$row]
}
